$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Current layout:  2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# Target layout:   2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The old "总计" sheet becomes the new "2022-Q1" fund-holding-detail
# sheet, and a brand new "总计" sheet is appended after it with the
# updated summary table (old summary rows + a new 2022-Q1 row).
# ------------------------------------------------------------------

$oldTotal = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

# Duplicate the whole "总计" sheet (this inherits sheetPr/pageMargins/
# column widths exactly) and place the duplicate right after it; this
# duplicate will become the new "总计" sheet.
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.ActiveSheet

# Rename the sheets into their final positions/names.
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$detail = $oldTotal

# ------------------------------------------------------------------
# Rebuild the "2022-Q1" sheet (formerly "总计") as a fund-holding
# detail sheet, matching the layout used by 2021-Q4, 2021-Q3, etc.
# ------------------------------------------------------------------
$detail.Cells.Clear()
$template.Range("A1:H7").Copy($detail.Range("A1:H7"))

$detail.Range("B1").Value = "基金代码"
$detail.Range("C1").Value = "基金名称"
$detail.Range("D1").Value = "基金规模"
$detail.Range("E1").Value = "股票总仓位"
$detail.Range("F1").Value = "仓位占比"
$detail.Range("G1").Value = "持有市值(亿元)"
$detail.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking codes/figures that must be
# stored as *text* (matching the source data), so they are entered with a
# leading apostrophe and then restored to the unstyled "Normal" look (the
# apostrophe entry alone would tag the cell with a quote-prefix style).
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$detail.Range("A2").Value = 0
Set-TextValue $detail.Range("B2") "004040"
$detail.Range("C2").Value = "金鹰医疗健康产业股票A"
Set-TextValue $detail.Range("D2") "10.91"
Set-TextValue $detail.Range("E2") "92.37"
Set-TextValue $detail.Range("F2") "4.83"
Set-TextValue $detail.Range("G2") "0.5270"
$detail.Range("H2").Value = 10

$detail.Range("A3").Value = 1
Set-TextValue $detail.Range("B3") "004041"
$detail.Range("C3").Value = "金鹰医疗健康产业股票C"
Set-TextValue $detail.Range("D3") "7.31"
Set-TextValue $detail.Range("E3") "92.37"
Set-TextValue $detail.Range("F3") "4.83"
Set-TextValue $detail.Range("G3") "0.3531"
$detail.Range("H3").Value = 10

$detail.Range("A4").Value = 2
Set-TextValue $detail.Range("B4") "012086"
$detail.Range("C4").Value = "博时健康生活混合型证券投资基金A"
Set-TextValue $detail.Range("D4") "3.52"
Set-TextValue $detail.Range("E4") "94.03"
Set-TextValue $detail.Range("F4") "2.70"
Set-TextValue $detail.Range("G4") "0.0950"
$detail.Range("H4").Value = 10

$detail.Range("A5").Value = 3
Set-TextValue $detail.Range("B5") "000523"
$detail.Range("C5").Value = "国投瑞银医疗保健行业灵活配置混合"
Set-TextValue $detail.Range("D5") "2.10"
Set-TextValue $detail.Range("E5") "85.78"
Set-TextValue $detail.Range("F5") "3.41"
Set-TextValue $detail.Range("G5") "0.0716"
$detail.Range("H5").Value = 10

$detail.Range("A6").Value = 4
Set-TextValue $detail.Range("B6") "005520"
$detail.Range("C6").Value = "国投瑞银创新医疗灵活配置混合"
Set-TextValue $detail.Range("D6") "0.51"
Set-TextValue $detail.Range("E6") "90.43"
Set-TextValue $detail.Range("F6") "3.98"
Set-TextValue $detail.Range("G6") "0.0203"
$detail.Range("H6").Value = 8

$detail.Range("A7").Value = 5
Set-TextValue $detail.Range("B7") "012087"
$detail.Range("C7").Value = "博时健康生活混合型证券投资基金C"
Set-TextValue $detail.Range("D7") "0.66"
Set-TextValue $detail.Range("E7") "94.03"
Set-TextValue $detail.Range("F7") "2.70"
Set-TextValue $detail.Range("G7") "0.0178"
$detail.Range("H7").Value = 10

# ------------------------------------------------------------------
# Rebuild the new "总计" sheet: same as the old one, but with a new
# 2022-Q1 row inserted at the top of the data (row 2), pushing the
# rest down by one row.
# ------------------------------------------------------------------
$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 6
$newTotal.Range("D2").Value = 1.08

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 13
$newTotal.Range("D3").Value = 17.18

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 3
$newTotal.Range("D4").Value = 2.49

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 24
$newTotal.Range("D5").Value = 37.71

$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 6
$newTotal.Range("D6").Value = 3.12

$newTotal.Range("A7").Value = 5
$newTotal.Range("B7").Value = "2020-Q4"
$newTotal.Range("C7").Value = 2
$newTotal.Range("D7").Value = 0
